$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes "SUBJECT - ORGANIZATION" section and
# everything below it down by one row) and populate it with the new
# "Citation Issued Location" field that documents
# ojb-cit-ext:CitationIssuedLocation.
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = "Citation Issued Loction"
$ws.Cells.Item(8, 2).Value = "Location where citation was issued."
$ws.Cells.Item(8, 5).Value = "cfm:CoreFilingMessage/ojb-cit-doc:CitationCase/ojb-cit-ext:CitationCaseAugmentation/ojb-cit-ext:Citation/ojb-cit-ext:CitationAugmentation/ojb-cit-ext:CitationIssuedLocation"

# Columns C and D are not applicable for this row, so make sure no stray
# formatting/content remains there (matches the other section-header-like
# rows in the sheet that simply omit unused cells).
$ws.Cells.Item(8, 3).Clear()
$ws.Cells.Item(8, 4).Clear()

# Match the row height used by the other rows in this two-line header band.
$ws.Rows.Item(8).RowHeight = 28

# Restore the selection/scroll position that the workbook was saved with.
[void]$ws.Range("C10").Select()
